$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'229.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'22.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.275"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05574"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.382"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.473"
$ws.Range("D7").Style = "Normal"
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = "'0.7827"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '7MXTokenMX'
$ws.Range("B9").Value = 'FTXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D9").Value = "'1.044"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '8FTXTokenFTT'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = "'0.1380"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = "'0.07396"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").Value = "'0.03163"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = "'0.02971"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = "'0.09278"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = "'0.001662"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '14BitForexTokenBF'
$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D16").Value = "'3.269"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '15MCDexMCB'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").Value = "'0.04788"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").Value = "'0.0005906"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '17OneONE'
$ws.Range("D19").Value = "'0.006241"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.005231"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.001064"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.0001501"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'3.930"
$ws.Range("D23").Style = "Normal"
$ws.Range("D26").Value = "'0.1244"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.0005005"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '26UpBotsUBXTBestin24h'
$ws.Range("D40").Value = "'0.04003"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.007022"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '40KickTokenKICK'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = "'0.1040"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = "'0.003327"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("D44").Value = "'0.009996"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005443"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.7860"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.04464"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.01011"
$ws.Range("D50").Style = "Normal"
